$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matches source data which uses
# literal text for prices, including thousand-separator dots).

$ws.Range("D2").Value = "58.958.85"
$ws.Range("E2").Value = "  -0.28%  "

$ws.Range("D3").Value = "2.499.57"
$ws.Range("E3").Value = "  +1.92%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.42"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.33"
$ws.Range("E6").Value = "  -2.76%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.571"
$ws.Range("E8").Value = "  +0.38%  "

$ws.Range("D9").Value = "2.521.69"
$ws.Range("E9").Value = "  +2.13%  "

$ws.Range("E10").Value = "  +1.07%  "

$ws.Range("E11").Value = "  +0.56%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.56"
$ws.Range("E12").Value = "  +5.81%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.353"
$ws.Range("E13").Value = "  +0.50%  "

$ws.Range("D14").Value = "2.940.63"
$ws.Range("E14").Value = "  +1.75%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.46"
$ws.Range("E15").Value = "  -2.44%  "

$ws.Range("D16").Value = "58.899.43"
$ws.Range("E16").Value = "  -0.18%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000139"
$ws.Range("E17").Value = "  +1.38%  "

$ws.Range("D18").Value = "2.521.20"
$ws.Range("E18").Value = "  -0.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.16"
$ws.Range("E19").Value = "  +0.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.26"
$ws.Range("E20").Value = "  -1.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.02"
$ws.Range("E21").Value = "  -0.32%  "

$ws.Range("E22").Value = "  +3.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.75"
$ws.Range("E23").Value = "  +0.47%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.76"
$ws.Range("E24").Value = "  +1.87%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.437"
$ws.Range("E25").Value = "  -4.81%  "

$ws.Range("E26").Value = "  +1.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  +1.98%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.85"
$ws.Range("E28").Value = "  +2.36%  "

$ws.Range("D29").Value = "0.0₃0771"
$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.80"
$ws.Range("E30").Value = "  -1.06%  "

$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.63"
$ws.Range("E31").Value = "  -1.12%  "

$ws.Range("E32").Value = "  -7.49%  "

$ws.Range("E33").Value = "  -0.10%  "

$ws.Range("E34").Value = "  +6.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "157.38"
$ws.Range("E35").Value = "  -0.57%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.61"
$ws.Range("E36").Value = "  +1.32%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.32"
$ws.Range("E37").Value = "  -4.55%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.59"
$ws.Range("E38").Value = "  -8.75%  "

$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.96"
$ws.Range("E39").Value = "  +1.35%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.62"
$ws.Range("E40").Value = "  -4.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "293.85"
$ws.Range("E41").Value = "  -7.17%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.68"
$ws.Range("E42").Value = "  -0.34%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.815"
$ws.Range("E43").Value = "  -2.40%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.995"
$ws.Range("E44").Value = "  -0.08%  "

$ws.Range("E45").Value = "  +2.19%  "

$ws.Range("E46").Value = "  +0.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0927"
$ws.Range("E47").Value = "  -1.26%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.14"
$ws.Range("E48").Value = "  +0.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.51"
$ws.Range("E49").Value = "  +0.16%  "

$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0227"
$ws.Range("E50").Value = "  -0.94%  "

$ws.Range("B51").Value = "Hedera"
$ws.Range("C51").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0512"
$ws.Range("E51").Value = "  -2.30%  "
